$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "211.43") need an
# explicit Text number format first, otherwise Excel auto-converts the
# assigned string into a numeric value instead of keeping it as text (the
# source data models these "Price" values as text, matching entries like
# "27.436.71" that cannot be parsed as a single number).
$forceTextCells = @("D5", "D8", "D15", "D16", "D18", "D20", "D22", "D25", "D31", "D37", "D38", "D39", "D40", "D43", "D49", "D50")
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.436.71"
$ws.Range("E2").Value = "  -1.88%  "
$ws.Range("D3").Value = "1.627.37"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").Value = "211.43"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("D8").Value = "22.98"
$ws.Range("E8").Value = "  -2.05%  "
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("E11").Value = "  -3.45%  "
$ws.Range("D12").Value = "1.857.95"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").Value = "1.627.99"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "0.557"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").Value = "64.94"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").Value = "27.443.99"
$ws.Range("E17").Value = "  -1.58%  "
$ws.Range("D18").Value = "228.64"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "0.0₃0718"
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("D20").Value = "7.52"
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Value = "10.77"
$ws.Range("E22").Value = "  +7.12%  "
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("E24").Value = "  +2.83%  "
$ws.Range("D25").Value = "149.48"
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("D31").Value = "0.0482"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("E32").Value = "  -1.08%  "
$ws.Range("D33").Value = "1.464.40"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("E35").Value = "  -1.31%  "
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("D37").Value = "0.559"
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "0.874"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0167"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").Value = "0.916"
$ws.Range("E40").Value = "  -1.13%  "
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("D43").Value = "67.76"
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("D47").Value = "1.767.42"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("D49").Value = "87.30"
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").Value = "0.0992"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("E51").Value = "  -2.01%  "
